# Auto-generated edit script applying the profit-table refresh from the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3040.9285
$ws.Range("I19").Value = 3252
$ws.Range("J19").Value = 2882.625
$ws.Range("K19").Value = 3252
$ws.Range("L19").Value = 2882.625
$ws.Range("M19").Value = -3077
$ws.Range("N19").Value = -3232.625
$ws.Range("H20").Value = 333338050
$ws.Range("J20").Value = 12024
$ws.Range("L20").Value = 12024
$ws.Range("N20").Value = -12484
$ws.Range("H28").Value = 959.0909
$ws.Range("I28").Value = 1265.2858
$ws.Range("J28").Value = 423.25
$ws.Range("K28").Value = 1265.2858
$ws.Range("L28").Value = 423.25
$ws.Range("M28").Value = -780.2858000000001
$ws.Range("N28").Value = -1393.25
$ws.Range("H35").Value = 333338050
$ws.Range("J35").Value = 12024
$ws.Range("L35").Value = 12024
$ws.Range("N35").Value = -12782
$ws.Range("H74").Value = 3916.3845
$ws.Range("I74").Value = 3357.5454
$ws.Range("J74").Value = 6990
$ws.Range("K74").Value = 3357.5454
$ws.Range("L74").Value = 6990
$ws.Range("M74").Value = -2421.5454
$ws.Range("N74").Value = -8862
$ws.Range("H77").Value = 3916.3845
$ws.Range("I77").Value = 3357.5454
$ws.Range("J77").Value = 6990
$ws.Range("K77").Value = 16787.727
$ws.Range("L77").Value = 34950
$ws.Range("M77").Value = -12107.727
$ws.Range("N77").Value = -44310
$ws.Range("H80").Value = 900.2778
$ws.Range("I80").Value = 712.6667
$ws.Range("K80").Value = 2138.0001
$ws.Range("M80").Value = -1140.0001
$ws.Range("H83").Value = 900.2778
$ws.Range("I83").Value = 712.6667
$ws.Range("K83").Value = 6414.0003
$ws.Range("M83").Value = -1422.0003
$ws.Range("H87").Value = 49462
$ws.Range("I87").Value = 15661
$ws.Range("K87").Value = 15661
$ws.Range("M87").Value = -14413
$ws.Range("H90").Value = 49462
$ws.Range("I90").Value = 15661
$ws.Range("K90").Value = 46983
$ws.Range("M90").Value = -40743
$ws.Range("H107").Value = 2685.1765
$ws.Range("I107").Value = 2674.7
$ws.Range("J107").Value = 2700.1428
$ws.Range("K107").Value = 2674.7
$ws.Range("L107").Value = 2700.1428
$ws.Range("M107").Value = -754.6999999999998
$ws.Range("N107").Value = -6540.1428
$ws.Range("H113").Value = 4037.75
$ws.Range("I113").Value = 2742
$ws.Range("K113").Value = 2742
$ws.Range("M113").Value = 512
$ws.Range("H132").Value = 4049.1667
$ws.Range("I132").Value = 4000.25
$ws.Range("K132").Value = 12000.75
$ws.Range("M132").Value = -9470.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1922.0769
$ws.Range("I2").Value = 1982.6666
$ws.Range("K2").Value = 1982.6666
$ws.Range("M2").Value = -1869.6666
$ws.Range("H45").Value = 2692.7
$ws.Range("I45").Value = 2686.4
$ws.Range("K45").Value = 2686.4
$ws.Range("M45").Value = -2309.4
$ws.Range("H61").Value = 7295.75
$ws.Range("J61").Value = 11191.5
$ws.Range("L61").Value = 11191.5
$ws.Range("N61").Value = -11615.5
$ws.Range("H63").Value = 2914.9
$ws.Range("I63").Value = 2594.0588
$ws.Range("J63").Value = 4733
$ws.Range("K63").Value = 2594.0588
$ws.Range("L63").Value = 4733
$ws.Range("M63").Value = -1908.0588
$ws.Range("N63").Value = -6105
$ws.Range("H66").Value = 2914.9
$ws.Range("I66").Value = 2594.0588
$ws.Range("J66").Value = 4733
$ws.Range("K66").Value = 12970.294
$ws.Range("L66").Value = 23665
$ws.Range("M66").Value = -9538.293999999998
$ws.Range("N66").Value = -30529
$ws.Range("H74").Value = 105381.266
$ws.Range("I74").Value = 117191
$ws.Range("K74").Value = 117191
$ws.Range("M74").Value = -116317
$ws.Range("H77").Value = 105381.266
$ws.Range("I77").Value = 117191
$ws.Range("K77").Value = 585955
$ws.Range("M77").Value = -581587
$ws.Range("H92").Value = 43275
$ws.Range("J92").Value = 36550
$ws.Range("L92").Value = 36550
$ws.Range("N92").Value = -41542
$ws.Range("H110").Value = 4185.115
$ws.Range("I110").Value = 4211.45
$ws.Range("J110").Value = 4097.3335
$ws.Range("K110").Value = 4211.45
$ws.Range("L110").Value = 4097.3335
$ws.Range("M110").Value = -2166.45
$ws.Range("N110").Value = -8187.3335
$ws.Range("H116").Value = 1922.0769
$ws.Range("I116").Value = 1982.6666
$ws.Range("K116").Value = 1982.6666
$ws.Range("M116").Value = 311.3334
$ws.Range("H122").Value = 2116.6775
$ws.Range("I122").Value = 2128.4614
$ws.Range("J122").Value = 2055.4
$ws.Range("K122").Value = 6385.3842
$ws.Range("L122").Value = 6166.200000000001
$ws.Range("M122").Value = -3935.3842
$ws.Range("N122").Value = -11066.2
$ws.Range("H132").Value = 2214.36
$ws.Range("I132").Value = 2067.9092
$ws.Range("J132").Value = 2329.4285
$ws.Range("K132").Value = 6203.7276
$ws.Range("L132").Value = 6988.2855
$ws.Range("M132").Value = -3673.7276
$ws.Range("N132").Value = -12048.2855
$ws.Range("H136").Value = 7295.75
$ws.Range("J136").Value = 11191.5
$ws.Range("L136").Value = 33574.5
$ws.Range("N136").Value = -38674.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1922.0769
$ws.Range("I3").Value = 1982.6666
$ws.Range("K3").Value = 1982.6666
$ws.Range("M3").Value = -1868.6666
$ws.Range("H94").Value = 505.2647
$ws.Range("I94").Value = 357.96667
$ws.Range("K94").Value = 357.96667
$ws.Range("M94").Value = 93.03332999999998
$ws.Range("H99").Value = 4271.645
$ws.Range("I99").Value = 3918.0625
$ws.Range("K99").Value = 3918.0625
$ws.Range("M99").Value = -2420.0625
$ws.Range("H105").Value = 1612.6666
$ws.Range("I105").Value = 1627.6364
$ws.Range("K105").Value = 1627.6364
$ws.Range("M105").Value = 119.3635999999999
$ws.Range("H107").Value = 2696.4
$ws.Range("I107").Value = 2677.2856
$ws.Range("K107").Value = 2677.2856
$ws.Range("M107").Value = -757.2856000000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 966
$ws.Range("I16").Value = 1207.6666
$ws.Range("K16").Value = 1207.6666
$ws.Range("M16").Value = -920.6666
$ws.Range("H31").Value = 386953.2
$ws.Range("I31").Value = 910810.75
$ws.Range("J31").Value = 2791
$ws.Range("K31").Value = 910810.75
$ws.Range("L31").Value = 2791
$ws.Range("M31").Value = -910515.75
$ws.Range("N31").Value = -3381
$ws.Range("H34").Value = 386953.2
$ws.Range("I34").Value = 910810.75
$ws.Range("J34").Value = 2791
$ws.Range("K34").Value = 910810.75
$ws.Range("L34").Value = 2791
$ws.Range("M34").Value = -910608.75
$ws.Range("N34").Value = -3195
$ws.Range("H58").Value = 2992.9395
$ws.Range("I58").Value = 2735.52
$ws.Range("J58").Value = 3797.375
$ws.Range("K58").Value = 2735.52
$ws.Range("L58").Value = 3797.375
$ws.Range("M58").Value = -2532.52
$ws.Range("N58").Value = -4203.375
$ws.Range("H62").Value = 3600
$ws.Range("J62").Value = 3600
$ws.Range("L62").Value = 3600
$ws.Range("N62").Value = -4848
$ws.Range("H65").Value = 3600
$ws.Range("J65").Value = 3600
$ws.Range("L65").Value = 18000
$ws.Range("N65").Value = -24240
$ws.Range("H99").Value = 380280.53
$ws.Range("I99").Value = 671677.8
$ws.Range("J99").Value = 16033.917
$ws.Range("K99").Value = 671677.8
$ws.Range("L99").Value = 16033.917
$ws.Range("M99").Value = -670179.8
$ws.Range("N99").Value = -19029.917
$ws.Range("H107").Value = 4584.528
$ws.Range("I107").Value = 720.6667
$ws.Range("J107").Value = 6516.4585
$ws.Range("K107").Value = 720.6667
$ws.Range("L107").Value = 6516.4585
$ws.Range("M107").Value = 1199.3333
$ws.Range("N107").Value = -10356.4585
$ws.Range("H113").Value = 966
$ws.Range("I113").Value = 1207.6666
$ws.Range("K113").Value = 1207.6666
$ws.Range("M113").Value = 962.3334
$ws.Range("H122").Value = 2969.0908
$ws.Range("I122").Value = 2969.0908
$ws.Range("K122").Value = 8907.2724
$ws.Range("M122").Value = -6457.2724
$ws.Range("H126").Value = 380280.53
$ws.Range("I126").Value = 671677.8
$ws.Range("J126").Value = 16033.917
$ws.Range("K126").Value = 2015033.4
$ws.Range("L126").Value = 48101.751
$ws.Range("M126").Value = -2012563.4
$ws.Range("N126").Value = -53041.751
$ws.Range("H132").Value = 2670.3572
$ws.Range("I132").Value = 2487.7778
$ws.Range("K132").Value = 7463.3334
$ws.Range("M132").Value = -4933.3334
$ws.Range("H136").Value = 2992.9395
$ws.Range("I136").Value = 2735.52
$ws.Range("J136").Value = 3797.375
$ws.Range("K136").Value = 8206.559999999999
$ws.Range("L136").Value = 11392.125
$ws.Range("M136").Value = -5656.559999999999
$ws.Range("N136").Value = -16492.125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 2200
$ws.Range("J88").Value = 2200
$ws.Range("L88").Value = 6600
$ws.Range("N88").Value = -7456
$ws.Range("H91").Value = 2200
$ws.Range("J91").Value = 2200
$ws.Range("L91").Value = 6600
$ws.Range("N91").Value = -9564
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 14630.143
$ws.Range("J99").Value = 9994.666999999999
$ws.Range("L99").Value = 29984.001
$ws.Range("N99").Value = -34476.001
$ws.Range("H122").Value = 1550.6
$ws.Range("I122").Value = 950
$ws.Range("J122").Value = 1700.75
$ws.Range("K122").Value = 8550
$ws.Range("L122").Value = 15306.75
$ws.Range("M122").Value = -6100
$ws.Range("N122").Value = -20206.75
$ws.Range("H131").Value = 12501437
$ws.Range("I131").Value = 100000936
$ws.Range("J131").Value = 1508.4572
$ws.Range("K131").Value = 300002808
$ws.Range("L131").Value = 4525.3716
$ws.Range("M131").Value = -299997768
$ws.Range("N131").Value = -14605.3716
$ws.Range("H132").Value = 6181.2085
$ws.Range("J132").Value = 1762.4546
$ws.Range("L132").Value = 15862.0914
$ws.Range("N132").Value = -20922.0914
$ws.Range("H137").Value = 862.4286
$ws.Range("I137").Value = 798.1667
$ws.Range("K137").Value = 2394.5001
$ws.Range("M137").Value = 2705.4999
$ws.Range("H139").Value = 3989.9167
$ws.Range("I139").Value = 2560.5454
$ws.Range("J139").Value = 5199.385
$ws.Range("K139").Value = 7681.6362
$ws.Range("L139").Value = 15598.155
$ws.Range("M139").Value = -2541.6362
$ws.Range("N139").Value = -25878.155
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 8006.5
$ws.Range("I7").Value = 8006.5
$ws.Range("K7").Value = 8006.5
$ws.Range("M7").Value = -7894.5
$ws.Range("H8").Value = 8006.5
$ws.Range("I8").Value = 8006.5
$ws.Range("K8").Value = 8006.5
$ws.Range("M8").Value = -7867.5
$ws.Range("H102").Value = 104970.1
$ws.Range("I102").Value = 5522.3335
$ws.Range("K102").Value = 5522.3335
$ws.Range("M102").Value = -3900.3335
$ws.Range("H113").Value = 3901
$ws.Range("I113").Value = 3861.4
$ws.Range("K113").Value = 3861.4
$ws.Range("M113").Value = -1691.4
$ws.Range("H122").Value = 6520.2974
$ws.Range("I122").Value = 4955.8125
$ws.Range("J122").Value = 7712.2856
$ws.Range("K122").Value = 14867.4375
$ws.Range("L122").Value = 23136.8568
$ws.Range("M122").Value = -12417.4375
$ws.Range("N122").Value = -28036.8568
$ws.Range("H126").Value = 9722.5
$ws.Range("I126").Value = 12466.25
$ws.Range("K126").Value = 37398.75
$ws.Range("M126").Value = -34928.75
$ws.Range("H132").Value = 42050.89
$ws.Range("I132").Value = 48573.176
$ws.Range("J132").Value = 4547.75
$ws.Range("K132").Value = 145719.528
$ws.Range("L132").Value = 13643.25
$ws.Range("M132").Value = -143189.528
$ws.Range("N132").Value = -18703.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1691.4468
$ws.Range("I16").Value = 1692.4117
$ws.Range("K16").Value = 1692.4117
$ws.Range("M16").Value = -1522.4117
$ws.Range("H35").Value = 2840.5715
$ws.Range("I35").Value = 2808.1667
$ws.Range("J35").Value = 3035
$ws.Range("K35").Value = 2808.1667
$ws.Range("L35").Value = 3035
$ws.Range("M35").Value = -2472.1667
$ws.Range("N35").Value = -3707
$ws.Range("H46").Value = 3252.5
$ws.Range("I46").Value = 2854.4
$ws.Range("J46").Value = 3916
$ws.Range("K46").Value = 2854.4
$ws.Range("L46").Value = 3916
$ws.Range("M46").Value = -2666.4
$ws.Range("N46").Value = -4292
$ws.Range("H61").Value = 2696.652
$ws.Range("I61").Value = 2680.8125
$ws.Range("J61").Value = 2732.8572
$ws.Range("K61").Value = 2680.8125
$ws.Range("L61").Value = 2732.8572
$ws.Range("M61").Value = -2478.8125
$ws.Range("N61").Value = -3136.8572
$ws.Range("H100").Value = 3992.3076
$ws.Range("I100").Value = 3991.6667
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 3991.6667
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -3450.6667
$ws.Range("N100").Value = -5082
$ws.Range("H111").Value = 74999
$ws.Range("J111").Value = 74999
$ws.Range("L111").Value = 74999
$ws.Range("N111").Value = -83179
$ws.Range("H113").Value = 2696.652
$ws.Range("I113").Value = 2680.8125
$ws.Range("J113").Value = 2732.8572
$ws.Range("K113").Value = 2680.8125
$ws.Range("L113").Value = 2732.8572
$ws.Range("M113").Value = -510.8125
$ws.Range("N113").Value = -7072.8572
$ws.Range("H136").Value = 3045.182
$ws.Range("I136").Value = 1314
$ws.Range("J136").Value = 6074.75
$ws.Range("K136").Value = 3942
$ws.Range("L136").Value = 18224.25
$ws.Range("M136").Value = -1392
$ws.Range("N136").Value = -23324.25
$ws.Range("H138").Value = 100429
$ws.Range("J138").Value = 100429
$ws.Range("L138").Value = 100429
$ws.Range("N138").Value = -110709
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2218.7222
$ws.Range("I81").Value = 2371.1875
$ws.Range("K81").Value = 4742.375
$ws.Range("M81").Value = -3681.375
$ws.Range("H84").Value = 2218.7222
$ws.Range("I84").Value = 2371.1875
$ws.Range("K84").Value = 23711.875
$ws.Range("M84").Value = -18407.875
$ws.Range("H100").Value = 1319.6471
$ws.Range("I100").Value = 1253.5
$ws.Range("K100").Value = 2507
$ws.Range("M100").Value = -1966
$ws.Range("H107").Value = 436
$ws.Range("J107").Value = 110
$ws.Range("L107").Value = 330
$ws.Range("N107").Value = -4170
$ws.Range("H122").Value = 2705.5557
$ws.Range("I122").Value = 2687.5625
$ws.Range("K122").Value = 8062.6875
$ws.Range("M122").Value = -5612.6875
$ws.Range("H126").Value = 2790.7917
$ws.Range("I126").Value = 2391.0527
$ws.Range("K126").Value = 7173.158100000001
$ws.Range("M126").Value = -4703.158100000001
$ws.Range("H132").Value = 968.25
$ws.Range("I132").Value = 967.9091
$ws.Range("K132").Value = 2903.7273
$ws.Range("M132").Value = -373.7273
$ws.Range("H136").Value = 457064.78
$ws.Range("I136").Value = 589904.2
$ws.Range("K136").Value = 1769712.6
$ws.Range("M136").Value = -1767162.6
